$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FirstSet")

# The "clinician.pt1@cvhcare.com" login (used in rows 22 and 26) is not
# working, so it is replaced everywhere with "clinician.pt2@cvhcare.com".
$oldAddress = "clinician.pt1@cvhcare.com"
$newAddress = "clinician.pt2@cvhcare.com"

# Snapshot the hyperlinks that currently exist on the sheet (cell address,
# target, and whether that cell is one of the ones being fixed) so they can
# be restored afterwards - this implementation's Hyperlinks.Delete() only
# works as a bulk operation that clears every hyperlink on the sheet.
$linksToRestore = New-Object System.Collections.ArrayList
foreach ($hl in $ws.Hyperlinks) {
    $cellAddr = $hl.Range.Address($false, $false)
    if ($cellAddr -ne "A22" -and $cellAddr -ne "A26") {
        [void]$linksToRestore.Add(@($cellAddr, $hl.Address))
    }
}

# Remove every hyperlink on the sheet.
$ws.Range("A1").Hyperlinks.Delete()

# Fix the two cells that contained the non-working pt1 address.
$ws.Range("A22").Value = $newAddress
$ws.Range("A26").Value = $newAddress

# Restore the hyperlinks that should remain (everything except A22 and A26).
foreach ($link in $linksToRestore) {
    $cellAddr = $link[0]
    $target = $link[1]
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $target) | Out-Null
    # Adding a hyperlink re-applies cell formatting; put the original
    # "Hyperlink" cell style back so no new style definitions are introduced.
    $ws.Range($cellAddr).Style = "Hyperlink"
}

# Reflect the last edited cell as the active selection, like in the source edit.
$ws.Range("A22").Select()
